$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row -> (DAMSLTag, DialogAct)
$updates = @{
    25  = @("ba", "Appreciation")
    30  = @("ba", "Appreciation")
    33  = @("ba", "Appreciation")
    55  = @("sv", "Statement-opinion")
    78  = @("sv", "Statement-opinion")
    93  = @("b", "Acknowledge (Backchannel)")
    94  = @("aa", "Agree/Accept")
    101 = @("%", "Uninterpretable")
    104 = @("sd", "Statement-non-opinion")
    108 = @("aa", "Agree/Accept")
    116 = @("sv", "Statement-opinion")
    117 = @("ba", "Appreciation")
    124 = @("sd", "Statement-non-opinion")
    137 = @("aa", "Agree/Accept")
    160 = @("sv", "Statement-opinion")
    163 = @("sd", "Statement-non-opinion")
    164 = @("aa", "Agree/Accept")
    166 = @("%", "Uninterpretable")
    184 = @("sd", "Statement-non-opinion")
    192 = @("sv", "Statement-opinion")
    194 = @("sv", "Statement-opinion")
    198 = @("sd", "Statement-non-opinion")
    220 = @("sd", "Statement-non-opinion")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
